# Auto-generated Excel COM-interop script
# Applies market-price / profit recalculation updates to the Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1040.9166
$ws.Range("I9").Value = 1094.2
$ws.Range("J9").Value = 774.5
$ws.Range("K9").Value = 1094.2
$ws.Range("L9").Value = 774.5
$ws.Range("M9").Value = -925.2
$ws.Range("N9").Value = -1112.5
$ws.Range("H19").Value = 2410
$ws.Range("I19").Value = 2410
$ws.Range("K19").Value = 2410
$ws.Range("M19").Value = -2235
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H80").Value = 1075
$ws.Range("I80").Value = 1075
$ws.Range("K80").Value = 3225
$ws.Range("M80").Value = -2227
$ws.Range("H83").Value = 1075
$ws.Range("I83").Value = 1075
$ws.Range("K83").Value = 9675
$ws.Range("M83").Value = -4683
$ws.Range("H86").Value = 4499.5
$ws.Range("I86").Value = 4999
$ws.Range("J86").Value = 4000
$ws.Range("K86").Value = 4999
$ws.Range("L86").Value = 4000
$ws.Range("M86").Value = -3876
$ws.Range("N86").Value = -6246
$ws.Range("H89").Value = 4499.5
$ws.Range("I89").Value = 4999
$ws.Range("J89").Value = 4000
$ws.Range("K89").Value = 24995
$ws.Range("L89").Value = 20000
$ws.Range("M89").Value = -19379
$ws.Range("N89").Value = -31232
$ws.Range("H106").Value = 10000
$ws.Range("I106").Value = 10000
$ws.Range("K106").Value = 10000
$ws.Range("M106").Value = -9369
$ws.Range("H107").Value = 324.22223
$ws.Range("I107").Value = 202.85715
$ws.Range("J107").Value = 749
$ws.Range("K107").Value = 202.85715
$ws.Range("L107").Value = 749
$ws.Range("M107").Value = 1717.14285
$ws.Range("N107").Value = -4589
$ws.Range("H127").Value = 797
$ws.Range("I127").Value = 797
$ws.Range("K127").Value = 2391
$ws.Range("M127").Value = 2569
$ws.Range("H134").Value = 124999
$ws.Range("J134").Value = 124999
$ws.Range("L134").Value = 124999
$ws.Range("N134").Value = -135139
$ws.Range("H137").Value = 2242.6667
$ws.Range("I137").Value = 1368.8462
$ws.Range("K137").Value = 4106.5386
$ws.Range("M137").Value = -1556.5386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3954.862
$ws.Range("I32").Value = 3953.2856
$ws.Range("K32").Value = 3953.2856
$ws.Range("M32").Value = -3666.2856
$ws.Range("H97").Value = 3222
$ws.Range("I97").Value = 2000
$ws.Range("K97").Value = 2000
$ws.Range("M97").Value = -1504
$ws.Range("H132").Value = 2988.36
$ws.Range("I132").Value = 2200.875
$ws.Range("J132").Value = 4388.3335
$ws.Range("K132").Value = 6602.625
$ws.Range("L132").Value = 13165.0005
$ws.Range("M132").Value = -4072.625
$ws.Range("N132").Value = -18225.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2099.4
$ws.Range("I20").Value = 1999.25
$ws.Range("K20").Value = 1999.25
$ws.Range("M20").Value = -1752.25
$ws.Range("H94").Value = 3100
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H105").Value = 4823
$ws.Range("I105").Value = 4823
$ws.Range("K105").Value = 4823
$ws.Range("M105").Value = -3076
$ws.Range("H107").Value = 1554.75
$ws.Range("I107").Value = 1370.3334
$ws.Range("K107").Value = 1370.3334
$ws.Range("M107").Value = 549.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2055.8572
$ws.Range("I31").Value = 2130.6667
$ws.Range("J31").Value = 1607
$ws.Range("K31").Value = 2130.6667
$ws.Range("L31").Value = 1607
$ws.Range("M31").Value = -1835.6667
$ws.Range("N31").Value = -2197
$ws.Range("H34").Value = 2055.8572
$ws.Range("I34").Value = 2130.6667
$ws.Range("J34").Value = 1607
$ws.Range("K34").Value = 2130.6667
$ws.Range("L34").Value = 1607
$ws.Range("M34").Value = -1928.6667
$ws.Range("N34").Value = -2011
$ws.Range("H43").Value = 20000
$ws.Range("J43").Value = 20000
$ws.Range("L43").Value = 20000
$ws.Range("N43").Value = -20368
$ws.Range("H58").Value = 2159.923
$ws.Range("I58").Value = 1690.8334
$ws.Range("J58").Value = 2562
$ws.Range("K58").Value = 1690.8334
$ws.Range("L58").Value = 2562
$ws.Range("M58").Value = -1487.8334
$ws.Range("N58").Value = -2968
$ws.Range("H86").Value = 7909.0713
$ws.Range("I86").Value = 9050.5
$ws.Range("J86").Value = 6387.1665
$ws.Range("K86").Value = 9050.5
$ws.Range("L86").Value = 6387.1665
$ws.Range("M86").Value = -7927.5
$ws.Range("N86").Value = -8633.166499999999
$ws.Range("H89").Value = 7909.0713
$ws.Range("I89").Value = 9050.5
$ws.Range("J89").Value = 6387.1665
$ws.Range("K89").Value = 45252.5
$ws.Range("L89").Value = 31935.8325
$ws.Range("M89").Value = -39636.5
$ws.Range("N89").Value = -43167.8325
$ws.Range("H96").Value = 60000
$ws.Range("J96").Value = 60000
$ws.Range("L96").Value = 60000
$ws.Range("N96").Value = -65492
$ws.Range("H97").Value = 46000
$ws.Range("J97").Value = 46000
$ws.Range("L97").Value = 46000
$ws.Range("N97").Value = -47982
$ws.Range("H101").Value = 20000
$ws.Range("J101").Value = 20000
$ws.Range("L101").Value = 20000
$ws.Range("N101").Value = -26490
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("H136").Value = 2159.923
$ws.Range("I136").Value = 1690.8334
$ws.Range("J136").Value = 2562
$ws.Range("K136").Value = 5072.5002
$ws.Range("L136").Value = 7686
$ws.Range("M136").Value = -2522.5002
$ws.Range("N136").Value = -12786

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 99
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 99
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 594
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -820
$ws.Range("H7").Value = 350
$ws.Range("J7").Value = 350
$ws.Range("L7").Value = 1050
$ws.Range("N7").Value = -1274
$ws.Range("H23").Value = 384
$ws.Range("J23").Value = 384
$ws.Range("L23").Value = 1152
$ws.Range("N23").Value = -1622
$ws.Range("H38").Value = 149.42857
$ws.Range("I38").Value = 50
$ws.Range("J38").Value = 166
$ws.Range("K38").Value = 150
$ws.Range("L38").Value = 498
$ws.Range("M38").Value = 197
$ws.Range("N38").Value = -1192
$ws.Range("H92").Value = 671.2
$ws.Range("I92").Value = 726.5
$ws.Range("K92").Value = 2179.5
$ws.Range("M92").Value = -931.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H122").Value = 1639
$ws.Range("I122").Value = 1674.5
$ws.Range("J122").Value = 1497
$ws.Range("K122").Value = 5023.5
$ws.Range("L122").Value = 4491
$ws.Range("M122").Value = -2573.5
$ws.Range("N122").Value = -9391

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2000
$ws.Range("I22").Value = 2000
$ws.Range("K22").Value = 2000
$ws.Range("M22").Value = -1705
$ws.Range("H27").Value = 2000
$ws.Range("I27").Value = 2000
$ws.Range("K27").Value = 2000
$ws.Range("M27").Value = -1893
$ws.Range("H40").Value = 3120.3845
$ws.Range("I40").Value = 3217.9167
$ws.Range("K40").Value = 3217.9167
$ws.Range("M40").Value = -3081.9167
$ws.Range("H55").Value = 297.8
$ws.Range("J55").Value = 297.8
$ws.Range("L55").Value = 297.8
$ws.Range("N55").Value = -643.8
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H136").Value = 3647.75
$ws.Range("I136").Value = 3647.75
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10943.25
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -8393.25
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7362.222
$ws.Range("I81").Value = 3608
$ws.Range("K81").Value = 7216
$ws.Range("M81").Value = -6155
$ws.Range("H84").Value = 7362.222
$ws.Range("I84").Value = 3608
$ws.Range("K84").Value = 36080
$ws.Range("M84").Value = -30776
$ws.Range("H132").Value = 2953.4
$ws.Range("I132").Value = 2276.5
$ws.Range("J132").Value = 3404.6667
$ws.Range("K132").Value = 6829.5
$ws.Range("L132").Value = 10214.0001
$ws.Range("M132").Value = -4299.5
$ws.Range("N132").Value = -15274.0001
$ws.Range("H136").Value = 2473.3333
$ws.Range("I136").Value = 2120
$ws.Range("J136").Value = 2826.6667
$ws.Range("K136").Value = 6360
$ws.Range("L136").Value = 8480.000100000001
$ws.Range("M136").Value = -3810
$ws.Range("N136").Value = -13580.0001
